$wb = $excel.ActiveWorkbook

# 1. Rename the "smile" header on the molecule sheet to "smiles" and move the
#    selection to the renamed cell (C1).
$wsMolecule = $wb.Worksheets.Item("molecule")
$wsMolecule.Activate()
$wsMolecule.Range("C1").Value = "smiles"
$wsMolecule.Range("C1").Select()

# 2. Select the header row on the COSY sheet (A1:K1) — matches the overlay
#    fix that highlights/selects the NMR-peak header columns.
$wsCosy = $wb.Worksheets.Item("COSY")
$wsCosy.Activate()
$wsCosy.Range("A1:K1").Select()

# 3. Add a new "NOESY" sheet after HMBC (the last existing sheet), matching
#    the layout/headers used by the other 2D-overlay sheets (COSY/HSQC/HMBC).
$wsHmbc = $wb.Worksheets.Item("HMBC")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsNoesy = $wb.Worksheets.Add($null, $lastSheet)
$wsNoesy.Name = "NOESY"

# Copy the header formatting (bold, border, centered) from HMBC's header row.
$wsHmbc.Range("B1:K1").Copy()
$wsNoesy.Range("B1:K1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill in the same column headers used by the other peak-list sheets.
$headers = @("f2 (ppm)", "f1 (ppm)", "Intensity", "Width f2", "Width f1", "Volume", "Type", "Flags", "Impurity/Compound", "Annotation")
$cols = @("B", "C", "D", "E", "F", "G", "H", "I", "J", "K")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $wsNoesy.Range($cols[$i] + "1").Value = $headers[$i]
}

# 4. Make NOESY the active sheet/tab with its header row selected.
$wsNoesy.Activate()
$wsNoesy.Range("A1:K1").Select()
